# Shift the "ser" blog reference numbers up by one: the post series
# previously ending at 178 now includes a new entry, 179.
#   H10: ser 176 -> 177
#   D10: ser 177 -> 178
#   B10: ser 178 -> 179  (new post)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H10").Value = "type: blog`nwidth: 2`nheight: 1`nser: 177"
$ws.Range("D10").Value = "type: blog`nwidth: 2`nheight: 1`nser: 178"
$ws.Range("B10").Value = "type: blog`nwidth: 2`nheight: 1`nser: 179"

# The author's cursor ends up on the newest entry, B10.
$ws.Range("B10").Select()
